$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("W4").Value = 1.44
# Row 5
$ws.Range("L5").Value = 1.57
$ws.Range("W5").Value = 1.27
# Row 6
$ws.Range("H6").Value = 2.66
$ws.Range("M6").Value = 1.12
$ws.Range("N6").Value = 2.6
$ws.Range("O6").Value = 1.51
$ws.Range("Q6").Value = 2.48
$ws.Range("T6").Value = 2.02
$ws.Range("U6").Value = 1.82
$ws.Range("X6").Value = 9.199999999999999
# Row 7
$ws.Range("I7").Value = 2.66
# Row 8
$ws.Range("G8").Value = 17.5
$ws.Range("K8").Value = 5.7
$ws.Range("N8").Value = 2.8
$ws.Range("P8").Value = 1.84
$ws.Range("Q8").Value = 1.98
$ws.Range("R8").Value = 1.32
$ws.Range("S8").Value = 3.5
$ws.Range("X8").Value = 970
$ws.Range("AL8").Value = 310
$ws.Range("AM8").Value = 370
# Row 9
$ws.Range("F9").Value = 1.31
$ws.Range("K9").Value = 13
$ws.Range("R9").Value = 1.63
$ws.Range("S9").Value = 1.93
# Row 11
$ws.Range("G11").Value = 7.2
$ws.Range("H11").Value = 1.67
$ws.Range("J11").Value = 3.5
$ws.Range("L11").Value = 1.01
$ws.Range("M11").Value = 1.09
$ws.Range("N11").Value = 2.94
$ws.Range("O11").Value = 1.43
$ws.Range("R11").Value = 1.23
$ws.Range("S11").Value = 3.9
$ws.Range("T11").Value = 2.12
$ws.Range("U11").Value = 1.74
$ws.Range("X11").Value = 11.5
$ws.Range("Y11").Value = 7
$ws.Range("Z11").Value = 9.800000000000001
$ws.Range("AA11").Value = 19.5
$ws.Range("AB11").Value = 17.5
$ws.Range("AC11").Value = 8.800000000000001
$ws.Range("AD11").Value = 11
$ws.Range("AE11").Value = 23
$ws.Range("AF11").Value = 55
$ws.Range("AG11").Value = 27
$ws.Range("AH11").Value = 970
$ws.Range("AI11").Value = 55
$ws.Range("AJ11").Value = 250
$ws.Range("AK11").Value = 140
$ws.Range("AL11").Value = 150
$ws.Range("AM11").Value = 240
$ws.Range("AO11").Value = 16.5
# Row 12
$ws.Range("G12").Value = 4.2
$ws.Range("L12").Value = 1.01
$ws.Range("M12").Value = 1.01
$ws.Range("N12").Value = 2.84
$ws.Range("O12").Value = 1.3
$ws.Range("P12").Value = 1.9
$ws.Range("R12").Value = 1.29
$ws.Range("S12").Value = 3.05
$ws.Range("T12").Value = 1.64
$ws.Range("U12").Value = 1.87
$ws.Range("V12").Value = 1.81
$ws.Range("W12").Value = 1.32
$ws.Range("X12").Value = 20
$ws.Range("Y12").Value = 13.5
$ws.Range("Z12").Value = 20
$ws.Range("AA12").Value = 38
$ws.Range("AB12").Value = 21
$ws.Range("AC12").Value = 11.5
$ws.Range("AD12").Value = 15.5
$ws.Range("AE12").Value = 34
$ws.Range("AF12").Value = 40
$ws.Range("AG12").Value = 23
$ws.Range("AH12").Value = 26
$ws.Range("AI12").Value = 55
$ws.Range("AJ12").Value = 100
$ws.Range("AK12").Value = 70
$ws.Range("AL12").Value = 85
$ws.Range("AM12").Value = 1000
$ws.Range("AN12").Value = 1000
$ws.Range("AO12").Value = 24
# Row 13
$ws.Range("I13").Value = 2.26
$ws.Range("J13").Value = 3.6
$ws.Range("L13").Value = 1.39
$ws.Range("M13").Value = 1.08
$ws.Range("N13").Value = 3.5
$ws.Range("O13").Value = 1.33
$ws.Range("P13").Value = 1.83
$ws.Range("Q13").Value = 2
$ws.Range("R13").Value = 1.32
$ws.Range("S13").Value = 3.5
$ws.Range("T13").Value = 1.77
$ws.Range("U13").Value = 2.04
$ws.Range("V13").Value = 1.8
$ws.Range("W13").Value = 1.37
$ws.Range("X13").Value = 13.5
$ws.Range("Z13").Value = 13.5
$ws.Range("AA13").Value = 29
$ws.Range("AC13").Value = 8
$ws.Range("AD13").Value = 11.5
$ws.Range("AE13").Value = 980
$ws.Range("AF13").Value = 26
$ws.Range("AG13").Value = 15.5
$ws.Range("AH13").Value = 18.5
$ws.Range("AI13").Value = 40
$ws.Range("AJ13").Value = 70
$ws.Range("AK13").Value = 44
$ws.Range("AL13").Value = 55
$ws.Range("AM13").Value = 110
$ws.Range("AN13").Value = 46
$ws.Range("AO13").Value = 19
# Row 14
$ws.Range("F14").Value = 4.6
$ws.Range("L14").Value = 1.4
$ws.Range("N14").Value = 3.55
$ws.Range("O14").Value = 1.33
$ws.Range("Q14").Value = 2
$ws.Range("R14").Value = 1.34
$ws.Range("S14").Value = 3.5
$ws.Range("V14").Value = 2.02
$ws.Range("W14").Value = 1.27
$ws.Range("X14").Value = 1000
$ws.Range("Y14").Value = 1000
$ws.Range("Z14").Value = 1000
$ws.Range("AA14").Value = 980
$ws.Range("AB14").Value = 1000
$ws.Range("AC14").Value = 970
$ws.Range("AD14").Value = 1000
$ws.Range("AE14").Value = 980
$ws.Range("AG14").Value = 1000
$ws.Range("AH14").Value = 1000
$ws.Range("AJ14").Value = 130
$ws.Range("AK14").Value = 70
$ws.Range("AL14").Value = 75
$ws.Range("AN14").Value = 75
$ws.Range("AO14").Value = 1000
# Row 15
$ws.Range("L15").Value = 1.36
$ws.Range("N15").Value = 3.95
$ws.Range("O15").Value = 1.27
$ws.Range("R15").Value = 1.4
$ws.Range("S15").Value = 3.05
$ws.Range("V15").Value = 1.19
$ws.Range("W15").Value = 2.42
$ws.Range("AM15").Value = 120
# Row 16
$ws.Range("F16").Value = 1.32
$ws.Range("G16").Value = 1.37
$ws.Range("H16").Value = 10.5
$ws.Range("J16").Value = 5.7
$ws.Range("K16").Value = 6.4
$ws.Range("L16").Value = 1.3
$ws.Range("N16").Value = 5.1
$ws.Range("P16").Value = 2.44
$ws.Range("V16").Value = 1.08
$ws.Range("W16").Value = 3.7
$ws.Range("Z16").Value = 120
$ws.Range("AA16").Value = 480
$ws.Range("AC16").Value = 15
$ws.Range("AD16").Value = 46
$ws.Range("AE16").Value = 200
# Row 17
$ws.Range("F17").Value = 1.97
$ws.Range("G17").Value = 2.74
$ws.Range("H17").Value = 3.2
$ws.Range("I17").Value = 5.5
$ws.Range("J17").Value = 2.66
$ws.Range("K17").Value = 5.6
$ws.Range("L17").Value = 1.01
$ws.Range("M17").Value = 1.01
$ws.Range("N17").Value = 2.06
$ws.Range("O17").Value = 1.28
$ws.Range("P17").Value = 1.33
$ws.Range("Q17").Value = 2.5
$ws.Range("R17").Value = 1.13
$ws.Range("S17").Value = 2.5
$ws.Range("T17").Value = 1.05
$ws.Range("U17").Value = 1.05
$ws.Range("V17").Value = 1.26
$ws.Range("W17").Value = 1.53
$ws.Range("X17").Value = 1000
$ws.Range("Y17").Value = 1000
$ws.Range("Z17").Value = 1000
$ws.Range("AA17").Value = 1000
$ws.Range("AB17").Value = 1000
$ws.Range("AC17").Value = 1000
$ws.Range("AD17").Value = 1000
$ws.Range("AE17").Value = 1000
$ws.Range("AF17").Value = 1000
$ws.Range("AG17").Value = 1000
$ws.Range("AH17").Value = 1000
$ws.Range("AI17").Value = 1000
$ws.Range("AJ17").Value = 1000
$ws.Range("AK17").Value = 1000
$ws.Range("AL17").Value = 1000
$ws.Range("AM17").Value = 1000
$ws.Range("AN17").Value = 1000
$ws.Range("AO17").Value = 1000
# Row 18
$ws.Range("L18").Value = 1.31
$ws.Range("R18").Value = 1.53
$ws.Range("S18").Value = 2.66
$ws.Range("U18").Value = 2.12
$ws.Range("V18").Value = 2.68
$ws.Range("W18").Value = 1.17
$ws.Range("AE18").Value = 18
$ws.Range("AN18").Value = 100
$ws.Range("AO18").Value = 7
# Row 20
$ws.Range("F20").Value = 1.44
$ws.Range("I20").Value = 11.5
# Row 22
$ws.Range("F22").Value = 2.22
$ws.Range("K22").Value = 3.45
# Row 23
$ws.Range("G23").Value = 3.3
